$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have a "Rule" / "Description" header pair (columns C/D).
# This upload renames the header row to "Name" / "Definition" so the table
# reads: ID, Tier, Name, Definition.
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Definition"

# Reflect the author's last on-screen selection/scroll position before saving.
$ws.Range("D3").Select()
